$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Produtos" (sheet1): drop the two extra product rows and the
# ReporEstoquepProd-only H column, reorder the price/method columns, and
# turn the remaining row into the new Ração/Pedigree record.
# ---------------------------------------------------------------------------
$wsProdutos = $wb.Worksheets.Item("Produtos")

$wsProdutos.Rows.Item(4).Delete()
$wsProdutos.Rows.Item(3).Delete()
$wsProdutos.Columns.Item(8).Delete()

$wsProdutos.Range("A1").Value = "Produto"
$wsProdutos.Range("B1").Value = "Marca"
$wsProdutos.Range("C1").Value = "Valor_Venda"
$wsProdutos.Range("D1").Value = "Valor_Compra"
$wsProdutos.Range("E1").Value = "Método_Venda"
$wsProdutos.Range("F1").Value = "Método_Compra"
$wsProdutos.Range("G1").Value = "ReporEstoquepProd"

$wsProdutos.Range("A2").Value = "Ração"
$wsProdutos.Range("B2").Value = "Pedigree"
$wsProdutos.Range("C2").Value = 20
$wsProdutos.Range("D2").Value = 10
$wsProdutos.Range("E2").Value = "Pacote"
$wsProdutos.Range("F2").Value = "Pacote"
$wsProdutos.Range("G2").Value = 10

# ---------------------------------------------------------------------------
# Sheet "Estoque" (sheet2): relabel the header columns and rewrite the
# single data row to match.
# ---------------------------------------------------------------------------
$wsEstoque = $wb.Worksheets.Item("Estoque")

$wsEstoque.Range("A1").Value = "Produto"
$wsEstoque.Range("B1").Value = "Marca"
$wsEstoque.Range("C1").Value = "Método"
$wsEstoque.Range("D1").Value = "Quantidade"
$wsEstoque.Range("E1").Value = "Valor_Compra"
$wsEstoque.Range("F1").Value = "Valor_Venda"

$wsEstoque.Range("A2").Value = "Ração"
$wsEstoque.Range("B2").Value = "Pedigree"
$wsEstoque.Range("C2").Value = "Pacote"
$wsEstoque.Range("D2").Value = 10
$wsEstoque.Range("E2").Value = 10
$wsEstoque.Range("F2").Value = 10

# ---------------------------------------------------------------------------
# Sheet "Vendas" (sheet3): replace the old header with the new sale-record
# columns and drop the now-unused 7th column.
# ---------------------------------------------------------------------------
$wsVendas = $wb.Worksheets.Item("Vendas")

$wsVendas.Range("A1").Value = "Num_Venda"
$wsVendas.Range("B1").Value = "Valor_Ganho"
$wsVendas.Range("C1").Value = "Frete"
$wsVendas.Range("D1").Value = "Desconto"
$wsVendas.Range("E1").Value = "Método_Pagamento"
$wsVendas.Range("F1").Value = "Comentário"
$wsVendas.Columns.Item(7).Delete()

# ---------------------------------------------------------------------------
# Sheet "Métodos" (sheet4): drop the extra "Pacote 3kg" row.
# ---------------------------------------------------------------------------
$wsMetodos = $wb.Worksheets.Item("Métodos")

$wsMetodos.Range("A1").Value = "Métodos"
$wsMetodos.Range("A2").Value = "Pacote"
$wsMetodos.Range("A3").Value = "Gramas"
$wsMetodos.Rows.Item(4).Delete()
